$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.067.76'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.567.42'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '208.87'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.14'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('E10').Value = '  +1.92%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0860'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '1.568.36'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.521'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').Value = '27.058.78'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '61.94'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '0.0₃0707'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.44'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '216.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.21'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '154.09'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.09'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0476'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('E30').Value = '  +4.84%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D33').Value = '1.425.65'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +13.27%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.62'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.35'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.83'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.814'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.33'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '64.93'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').Value = '1.704.08'
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '86.75'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('E48').Value = '  +3.48%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0519'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0965'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('E51').Value = '  +0.48%  '
